# Implements the "refactor: Uses EST# to match invoices to citibuy and renames column" change
# on the SampleAgingReport.xlsx test fixture.
#
# Summary of the data change (per the diff):
#   - Column A/B ("EST#"/"Invoice") now hold the plain invoice numbers (Invoice#2, Invoice#3,
#     Invoice#6, #1, #2, #4, #5) instead of the old Vendor-ID-derived strings.
#   - Column C ("WO") now holds plain numeric work-order numbers (101,102,103,201,202,203,204)
#     instead of the old "######-R" strings.
#   - Column I ("PO: Release") values shift because some of the now-unused shared strings were
#     removed.
#   - Row 19 (a blank trailing row) is removed, so the sheet now spans A1:I18.
#   - The active selection moves from I5 to B9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows 2-8 -------------------------------------------------

# Row 2
$ws.Range("A2").Value = "Invoice#2"
$ws.Range("B2").Value = "Invoice#2"
$ws.Range("C2").Value = 101
$ws.Range("I2").Value = "P111:1"

# Row 3
$ws.Range("A3").Value = "Invoice#3"
$ws.Range("B3").Value = "Invoice#3"
$ws.Range("C3").Value = 102
$ws.Range("I3").Value = "P111:1"

# Row 4
$ws.Range("A4").Value = "Invoice#6"
$ws.Range("B4").Value = "Invoice#6"
$ws.Range("C4").Value = 103
$ws.Range("I4").Value = "P999"

# Row 5
$ws.Range("A5").Value = "#1"
$ws.Range("B5").Value = "#1"
$ws.Range("C5").Value = 201
$ws.Range("I5").Value = "P222"

# Row 6
$ws.Range("A6").Value = "#2"
$ws.Range("B6").Value = "#2"
$ws.Range("C6").Value = 202
$ws.Range("I6").Value = "P222"

# Row 7
$ws.Range("A7").Value = "#4"
$ws.Range("B7").Value = "#4"
$ws.Range("C7").Value = 203
$ws.Range("I7").Value = "P444:1"

# Row 8
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 204
$ws.Range("I8").Value = "P444:1"

# --- Remove the now-empty trailing row 19 ---------------------------------

$ws.Rows.Item(19).Delete()

# --- Update the saved selection -------------------------------------------

$ws.Range("B9").Select()
